$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 780, shifting the existing 2026/12/29.. rows (and everything
# below) down by one. This makes room for a new 2026/02/06 data point that sits
# right after the two existing 2026/02/06 rows (778, 779).
$ws.Rows.Item(780).Insert()

# Populate the newly inserted row with the new reading for 2026/02/06 (Friday).
# Force A780 to be stored as text (not auto-converted to a date serial number)
# so it matches the other date cells in the column, then restore the default
# "Normal" style so no stray number-format style is left behind.
$ws.Range("A780").NumberFormat = "@"
$ws.Range("A780").Value = "2026/02/06"
$ws.Range("A780").Style = "Normal"

$ws.Range("B780").Value = "金"
$ws.Range("C780").Value = 14
$ws.Range("D780").Value = 201
